$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set D6, D7, D8, D9 to 1 with an integer ("0") number format, ---
# --- keeping their existing borders/alignment (style dedup produces ---
# --- the two new cellXfs entries expected by the diff). ---
$ws.Range("D6").Value = 1
$ws.Range("D6").NumberFormat = "0"

$ws.Range("D7").Value = 1
$ws.Range("D7").NumberFormat = "0"

$ws.Range("D8").Value = 1
$ws.Range("D8").NumberFormat = "0"

$ws.Range("D9").Value = 1
$ws.Range("D9").NumberFormat = "0"

# --- Zoom out the view and move the frozen pane's visible column back ---
# --- to column B (first column after the column-A freeze). ---
$excel.ActiveWindow.Zoom = 70

# --- Drop the color-scale conditional formatting that used to cover ---
# --- E6:K9 and M7:S9 (the other colorScale ranges are left as-is). ---
$ws.Range("E6:K9").FormatConditions.Delete()
$ws.Range("M7:S9").FormatConditions.Delete()
